$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new quarterly data row (row 71) for period 01-04-2021
$ws.Range("A71").Value = "'01-04-2021"
$ws.Range("A71").Style = "Normal"

$ws.Range("B71").Value = 41289
$ws.Range("C71").Value = 20181
$ws.Range("D71").Value = 16683
$ws.Range("E71").Value = 16
$ws.Range("F71").Value = 457
$ws.Range("G71").Value = 2840
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 586
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 4
$ws.Range("L71").Value = 522

Write-Output "Row 71 written"
